$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 80 ("Primera"/"Segunda" for
# the previous week, 2021-11-29 / serial 44529). This pushes the old
# rows 80-81 down to 82-83 and opens up rows 80-81 for the new week's
# data (2022-01-24 / serial 44585).
$ws.Rows("80:81").Insert()

# New row 80: "Primera" quality for the new week, replacing the old
# week's row that used to sit at row 80.
$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(80, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value = "Metropolitana"
$ws.Cells.Item(80, 4).Value = 44585
$ws.Cells.Item(80, 5).Value = 13
$ws.Cells.Item(80, 6).Value = 100114007
$ws.Cells.Item(80, 7).Value = "Jengibre"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 790
$ws.Cells.Item(80, 11).Value = 10000
$ws.Cells.Item(80, 12).Value = 11000
$ws.Cells.Item(80, 13).Value = 10494
$ws.Cells.Item(80, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(80, 15).Value = "Perú"
$ws.Cells.Item(80, 16).Value = 807
$ws.Cells.Item(80, 17).Value = 13
$ws.Cells.Item(80, 18).Value = "Hortaliza"

# New row 81: "Segunda" quality for the new week.
$ws.Cells.Item(81, 1).Value = 9
$ws.Cells.Item(81, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44585
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100114007
$ws.Cells.Item(81, 7).Value = "Jengibre"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Segunda"
$ws.Cells.Item(81, 10).Value = 340
$ws.Cells.Item(81, 11).Value = 9000
$ws.Cells.Item(81, 12).Value = 9000
$ws.Cells.Item(81, 13).Value = 9000
$ws.Cells.Item(81, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(81, 15).Value = "Perú"
$ws.Cells.Item(81, 16).Value = 692
$ws.Cells.Item(81, 17).Value = 13
$ws.Cells.Item(81, 18).Value = "Hortaliza"
